$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10: Dương Tử Phúc advice text updated
$ws.Range("B10").Value = "Bạn là người thông minh tuy nhiên luôn đề cao cái tôi và luôn cố gắng thắng trong mọi cuộc chơi nên dễ gây ra nhiều thị phi. Bạn nên đề cao sự tử tế và nhu thuận hơn thì sẽ gặp nhiều sự hỗ trợ từ mọi người."

# Row 11: Tang Tuế Điếu advice text updated
$ws.Range("B11").Value = "Bạn là người năng động luôn cố gắng xoay sở để có thể cải thiện vị thế bản thân, tuy nhiên cần cẩn trọng việc đi lại, giữ hòa khí với mọi người."

# Row 9: Âm Long Trực advice text updated + new C column value
$ws.Range("B9").Value = "Bạn là người thông minh, biết cách ứng xử phù hợp và nên giữ đức tính nhu thuận làm kim chỉ nam cuộc đời để gặp nhiều may mắn. Thuận thiên vô chiến tự nhiên thành."
$ws.Range("C9").Value = "Đặc biệt điều tiên quyét là bạn phải biết chọn bạn mà chơi chọn thầy mà theo thì bạn mới gặt được thành quả tốt đẹp."

# Row 8: Tuế Hổ Phù advice text updated
$ws.Range("B8").Value = "Bạn sinh ra gặp rất nhiều sóng gió cuộc đời nhưng đừng nản chí và bỏ cuộc vì đến khi vào đại vận bạn sẽ đủ chín chắn, trưởng thành do tích lũy trước đó. Bạn sẽ được hưởng trọn vẹn thành quả của sự cố gắng, cuộc sống gắn liền phần nhiều đến tín ngưỡng và tôn giáo. Bạn hãy cố gắng trau dồi kinh nghiệm, kiến thức và trải nghiệm."

$ws.Range("J26").Select()

$wb.Save()
